# #5: cash & deposit done
# Rebuild the "存款" (deposit) sheet: add bank / deposit_type / currency
# headers plus the common metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) that every
# other sheet in this workbook already carries, and drop the old
# (mostly-empty) "quantity" column in favour of the totals that used to
# sit in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# ---- header row -----------------------------------------------------
$headers = @{
    B = "bank"
    C = "deposit_type"
    D = "currency"
    E = "owner"
    F = "total"
    G = "property_category"
    H = "category"
    I = "date"
    J = "legislator_name"
    K = "legislator_id"
    L = "source_file"
    M = "index"
}

$headerStyleSrc = $ws.Range("B1")

foreach ($col in @("B","C","D","E","F","G","H","I","J","K","L","M")) {
    $cell = $ws.Range($col + "1")
    $cell.Value = $headers[$col]
    $headerStyleSrc.Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null
}

# ---- data rows --------------------------------------------------------
$rows = @(
    @{ row=2; A=89; B="安泰商業銀行";         C="活期儲蓄存款"; D="新臺幣"; F=6098718;    M=89 }
    @{ row=3; A=90; B="安泰商業銀行";         C="支票存款";     D="新臺幣"; F=161;        M=90 }
    @{ row=4; A=91; B="華泰商業銀行";         C="活期儲蓄存款"; D="新臺幣"; F=29425;      M=91 }
    @{ row=5; A=92; B="華泰商業銀行";         C="支票存款";     D="新臺幣"; F=59831;      M=92 }
    @{ row=6; A=93; B="中國銀行";             C="綜合存款";     D="人民幣"; F=7879500;    M=93 }
    @{ row=7; A=94; B="美商美國銀行";         C="綜合存款";     D="美金";   F=30533705.8; M=94 }
    @{ row=8; A=95; B="永豐商業銀行敦南分行"; C="活期儲蓄存款"; D="新臺幣"; F=1593775;    M=95 }
)

$indexStyleSrc = $ws.Range("A2")

foreach ($r in $rows) {
    $rn = $r.row

    $aCell = $ws.Cells.Item($rn, 1)
    $aCell.Value = $r.A
    $indexStyleSrc.Copy() | Out-Null
    $aCell.PasteSpecial(-4122) | Out-Null

    $ws.Range("B" + $rn).Value = $r.B
    $ws.Range("C" + $rn).Value = $r.C
    $ws.Range("D" + $rn).Value = $r.D
    $ws.Range("E" + $rn).Value = "羅淑蕾"
    $ws.Range("F" + $rn).Value = $r.F
    $ws.Range("G" + $rn).Value = "deposit"
    $ws.Range("H" + $rn).Value = "normal"

    # Force text so Excel doesn't reinterpret the ISO date string as a
    # serial date number.
    $dateCell = $ws.Range("I" + $rn)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2011-11-21"

    $ws.Range("J" + $rn).Value = "羅淑蕾"
    $ws.Range("K" + $rn).Value = 1638
    $ws.Range("L" + $rn).Value = "tmpe5cc1"
    $ws.Range("M" + $rn).Value = $r.M
}
